$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Two new data rows (168, 169) continuing the CRM180 batch series that
#    ends at row 167. Copy direct formatting from row 167 first (reuses the
#    existing style indices s=5 / s=1 instead of minting new ones), then
#    fill in the values/formula for the new rows.
# ---------------------------------------------------------------------------
$ws.Range("A167:F167").Copy() | Out-Null
$ws.Range("A168:F169").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A168").Value = 20250819
$ws.Range("B168").Value = 2219.70714
$ws.Range("C168").Value = 2224.4699999999998
$ws.Range("E168").Value = 180
$ws.Range("F168").Value = "CRM180_opened20250714_PP"

$ws.Range("A169").Value = 20250819
$ws.Range("B169").Value = 2213.6607899999999
$ws.Range("C169").Value = 2224.4699999999998
$ws.Range("E169").Value = 181
$ws.Range("F169").Value = "CRM180_opened20250714_PP"

# Setting the whole D168:D169 block in one shot lets the engine recognise the
# repeated formula pattern and store it as a shared formula, matching the
# existing shared group that already covers D162:D167.
$ws.Range("D168:D169").Formula = "=100*(B168-C168)/C168"

$ws.Range("A168:F169").RowHeight = 15.6

# ---------------------------------------------------------------------------
# 2. Three small formatted-but-empty rows further down (172-174), using a
#    new "Lucida Console" 8pt font - the second of which also carries a
#    solid white fill.
# ---------------------------------------------------------------------------
$ws.Range("D172").Font.Name = "Lucida Console"
$ws.Range("D172").Font.Size = 8
$ws.Range("D172").Font.Color = 0
$ws.Range("D172").VerticalAlignment = -4108   # xlCenter

$ws.Range("D173").Font.Name = "Lucida Console"
$ws.Range("D173").Font.Size = 8
$ws.Range("D173").Font.Color = 0
$ws.Range("D173").VerticalAlignment = -4108   # xlCenter

$ws.Range("D174").Font.Name = "Lucida Console"
$ws.Range("D174").Font.Size = 8
$ws.Range("D174").Font.Color = 0
$ws.Range("D174").Interior.Color = 16777215   # white
$ws.Range("D174").VerticalAlignment = -4108   # xlCenter

# ---------------------------------------------------------------------------
# 3. Update the view: scroll position & active selection.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 154
$win.ScrollColumn = 1
$ws.Range("M170").Select() | Out-Null
